$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 106 (existing rows 106-124 shift down to 108-126),
# mirroring the weekly Kiwi price data that was added for market date 44476.
$ws.Range("A106:A107").EntireRow.Insert()

# Row 106: "Especial" quality entry for the new date
$ws.Range("A106").Value = 4
$ws.Range("B106").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C106").Value = "Los Lagos"
$ws.Range("D106").Value = 44476
$ws.Range("E106").Value = 10
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100101
$ws.Range("H106").Value = "Berries"
$ws.Range("I106").Value = 100101007
$ws.Range("J106").Value = "Kiwi"
$ws.Range("K106").Value = "Hayward"
$ws.Range("L106").Value = "Especial"
$ws.Range("M106").Value = 200
$ws.Range("N106").Value = 21000
$ws.Range("O106").Value = 21000
$ws.Range("P106").Value = 21000
$ws.Range("Q106").Value = "$/caja 15 kilos"
$ws.Range("R106").Value = "Provincia de Curicó"
$ws.Range("S106").Value = 1400
$ws.Range("T106").Value = 15

# Row 107: "Primera" quality entry for the same new date
$ws.Range("A107").Value = 4
$ws.Range("B107").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C107").Value = "Los Lagos"
$ws.Range("D107").Value = 44476
$ws.Range("E107").Value = 10
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100101
$ws.Range("H107").Value = "Berries"
$ws.Range("I107").Value = 100101007
$ws.Range("J107").Value = "Kiwi"
$ws.Range("K107").Value = "Hayward"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 300
$ws.Range("N107").Value = 15000
$ws.Range("O107").Value = 16000
$ws.Range("P107").Value = 15500
$ws.Range("Q107").Value = "$/caja 15 kilos"
$ws.Range("R107").Value = "Provincia de Curicó"
$ws.Range("S107").Value = 1033
$ws.Range("T107").Value = 15
